# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for rows 2-44 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 2
    8  = 3
    9  = 0
    10 = 2
    11 = 1
    12 = 3
    13 = 6
    14 = 1
    15 = 7
    16 = 4
    17 = 3
    18 = 1
    19 = 4
    20 = 11
    21 = 3
    22 = 3
    23 = 4
    24 = 4
    25 = 8
    26 = 5
    27 = 6
    28 = 5
    29 = 5
    30 = 5
    31 = 3
    32 = 5
    33 = 9
    34 = 7
    35 = 7
    36 = 6
    37 = 4
    38 = 4
    39 = 10
    40 = 2
    41 = 3
    42 = 5
    43 = 4
    44 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
